$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the Espárragos weekly block (row 71),
# pushing the existing rows 71-94 down to 74-97.
$ws.Range("A71:R73").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new week's data
# (constants shared by every row in this data block)
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$fecha = 44511
$codreg = 13
$catId = 300000000
$categoria = "Espárragos"
$variedad = "Sin especificar"
$unidad = "`$/kilo"
$kgUnidades = 1
$clasificacion = "Hortaliza"

# Row 71: Banquete
$ws.Range("A71").Value = $mercadoId
$ws.Range("B71").Value = $mercado
$ws.Range("C71").Value = $region
$ws.Range("D71").Value = $fecha
$ws.Range("E71").Value = $codreg
$ws.Range("F71").Value = $catId
$ws.Range("G71").Value = $categoria
$ws.Range("H71").Value = $variedad
$ws.Range("I71").Value = "Banquete"
$ws.Range("J71").Value = 460
$ws.Range("K71").Value = 1200
$ws.Range("L71").Value = 1300
$ws.Range("M71").Value = 1250
$ws.Range("N71").Value = $unidad
$ws.Range("O71").Value = "Provincia de Linares"
$ws.Range("P71").Value = 1250
$ws.Range("Q71").Value = $kgUnidades
$ws.Range("R71").Value = $clasificacion

# Row 72: Primera
$ws.Range("A72").Value = $mercadoId
$ws.Range("B72").Value = $mercado
$ws.Range("C72").Value = $region
$ws.Range("D72").Value = $fecha
$ws.Range("E72").Value = $codreg
$ws.Range("F72").Value = $catId
$ws.Range("G72").Value = $categoria
$ws.Range("H72").Value = $variedad
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 340
$ws.Range("K72").Value = 1000
$ws.Range("L72").Value = 1100
$ws.Range("M72").Value = 1050
$ws.Range("N72").Value = $unidad
$ws.Range("O72").Value = "Provincia de Linares"
$ws.Range("P72").Value = 1050
$ws.Range("Q72").Value = $kgUnidades
$ws.Range("R72").Value = $clasificacion

# Row 73: Segunda
$ws.Range("A73").Value = $mercadoId
$ws.Range("B73").Value = $mercado
$ws.Range("C73").Value = $region
$ws.Range("D73").Value = $fecha
$ws.Range("E73").Value = $codreg
$ws.Range("F73").Value = $catId
$ws.Range("G73").Value = $categoria
$ws.Range("H73").Value = $variedad
$ws.Range("I73").Value = "Segunda"
$ws.Range("J73").Value = 250
$ws.Range("K73").Value = 800
$ws.Range("L73").Value = 900
$ws.Range("M73").Value = 840
$ws.Range("N73").Value = $unidad
$ws.Range("O73").Value = "Provincia de Linares"
$ws.Range("P73").Value = 840
$ws.Range("Q73").Value = $kgUnidades
$ws.Range("R73").Value = $clasificacion
